$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1659.7222
$ws.Range("I38").Value = 794.1667
$ws.Range("J38").Value = 2092.5
$ws.Range("K38").Value = 2382.5001
$ws.Range("L38").Value = 6277.5
$ws.Range("M38").Value = -2010.5001
$ws.Range("N38").Value = -7021.5
$ws.Range("H118").Value = 381.2857
$ws.Range("I118").Value = 381.2857
$ws.Range("K118").Value = 1143.8571
$ws.Range("M118").Value = 513.1428999999998
$ws.Range("H132").Value = 6240.222
$ws.Range("I132").Value = 5123.143
$ws.Range("J132").Value = 10150
$ws.Range("K132").Value = 15369.429
$ws.Range("L132").Value = 30450
$ws.Range("M132").Value = -12839.429
$ws.Range("N132").Value = -35510
$ws.Range("H134").Value = 111183416
$ws.Range("J134").Value = 111183416
$ws.Range("L134").Value = 111183416
$ws.Range("N134").Value = -111193556
$ws.Range("H138").Value = 1927.7028
$ws.Range("I138").Value = 1521.3055
$ws.Range("J138").Value = 2312.7104
$ws.Range("K138").Value = 4563.916499999999
$ws.Range("L138").Value = 6938.1312
$ws.Range("M138").Value = 576.0835000000006
$ws.Range("N138").Value = -17218.1312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1694.0204
$ws.Range("I61").Value = 1571.6571
$ws.Range("J61").Value = 1999.9286
$ws.Range("K61").Value = 1571.6571
$ws.Range("L61").Value = 1999.9286
$ws.Range("M61").Value = -1359.6571
$ws.Range("N61").Value = -2423.9286
$ws.Range("H74").Value = 1834.807
$ws.Range("I74").Value = 1853.9215
$ws.Range("J74").Value = 1672.3334
$ws.Range("K74").Value = 1853.9215
$ws.Range("L74").Value = 1672.3334
$ws.Range("M74").Value = -979.9214999999999
$ws.Range("N74").Value = -3420.3334
$ws.Range("H77").Value = 1834.807
$ws.Range("I77").Value = 1853.9215
$ws.Range("J77").Value = 1672.3334
$ws.Range("K77").Value = 9269.6075
$ws.Range("L77").Value = 8361.666999999999
$ws.Range("M77").Value = -4901.6075
$ws.Range("N77").Value = -17097.667
$ws.Range("H101").Value = 29900
$ws.Range("J101").Value = 29900
$ws.Range("L101").Value = 29900
$ws.Range("N101").Value = -36390
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524
$ws.Range("H109").Value = 17144.857
$ws.Range("J109").Value = 17144.857
$ws.Range("L109").Value = 17144.857
$ws.Range("N109").Value = -19918.857
$ws.Range("H112").Value = 19161.666
$ws.Range("J112").Value = 19161.666
$ws.Range("L112").Value = 19161.666
$ws.Range("N112").Value = -22115.666
$ws.Range("H114").Value = 38695
$ws.Range("J114").Value = 38695
$ws.Range("L114").Value = 38695
$ws.Range("N114").Value = -47373
$ws.Range("H132").Value = 747984
$ws.Range("I132").Value = 1223435.1
$ws.Range("J132").Value = 6280.4
$ws.Range("K132").Value = 3670305.3
$ws.Range("L132").Value = 18841.2
$ws.Range("M132").Value = -3667775.3
$ws.Range("N132").Value = -23901.2
$ws.Range("H136").Value = 1694.0204
$ws.Range("I136").Value = 1571.6571
$ws.Range("J136").Value = 1999.9286
$ws.Range("K136").Value = 4714.971299999999
$ws.Range("L136").Value = 5999.7858
$ws.Range("M136").Value = -2164.971299999999
$ws.Range("N136").Value = -11099.7858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 427
$ws.Range("I36").Value = 427
$ws.Range("K36").Value = 427
$ws.Range("M36").Value = 107
$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344
$ws.Range("H112").Value = 23000
$ws.Range("J112").Value = 23000
$ws.Range("L112").Value = 23000
$ws.Range("N112").Value = -25954
$ws.Range("H134").Value = 3553.377
$ws.Range("I134").Value = 942.8
$ws.Range("J134").Value = 7067.615
$ws.Range("K134").Value = 2828.4
$ws.Range("L134").Value = 21202.845
$ws.Range("M134").Value = -293.3999999999996
$ws.Range("N134").Value = -26272.845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 16817.75
$ws.Range("J43").Value = 16817.75
$ws.Range("L43").Value = 16817.75
$ws.Range("N43").Value = -17185.75
$ws.Range("H101").Value = 16817.75
$ws.Range("J101").Value = 16817.75
$ws.Range("L101").Value = 16817.75
$ws.Range("N101").Value = -23307.75
$ws.Range("H132").Value = 2581.4634
$ws.Range("I132").Value = 2001.8096
$ws.Range("J132").Value = 3190.1
$ws.Range("K132").Value = 6005.4288
$ws.Range("L132").Value = 9570.299999999999
$ws.Range("M132").Value = -3475.4288
$ws.Range("N132").Value = -14630.3
$ws.Range("H134").Value = 2103.7334
$ws.Range("I134").Value = 1243.7894
$ws.Range("J134").Value = 3589.0908
$ws.Range("K134").Value = 3731.3682
$ws.Range("L134").Value = 10767.2724
$ws.Range("M134").Value = -1196.3682
$ws.Range("N134").Value = -15837.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3914.913
$ws.Range("I137").Value = 1571
$ws.Range("J137").Value = 5717.923
$ws.Range("K137").Value = 4713
$ws.Range("L137").Value = 17153.769
$ws.Range("M137").Value = 387
$ws.Range("N137").Value = -27353.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 6746.1665
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H104").Value = 32215.25
$ws.Range("J104").Value = 32215.25
$ws.Range("L104").Value = 32215.25
$ws.Range("N104").Value = -39203.25
$ws.Range("H105").Value = 32000
$ws.Range("J105").Value = 32000
$ws.Range("L105").Value = 32000
$ws.Range("N105").Value = -38988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 909813.6
$ws.Range("I46").Value = 416.66666
$ws.Range("J46").Value = 2001090
$ws.Range("K46").Value = 416.66666
$ws.Range("L46").Value = 2001090
$ws.Range("M46").Value = -228.66666
$ws.Range("N46").Value = -2001466
$ws.Range("H101").Value = 13995.25
$ws.Range("J101").Value = 13995.25
$ws.Range("L101").Value = 13995.25
$ws.Range("N101").Value = -20485.25
$ws.Range("H105").Value = 32807.5
$ws.Range("J105").Value = 32807.5
$ws.Range("L105").Value = 32807.5
$ws.Range("N105").Value = -39795.5
$ws.Range("H110").Value = 22563.625
$ws.Range("J110").Value = 22563.625
$ws.Range("L110").Value = 22563.625
$ws.Range("N110").Value = -30743.625
$ws.Range("H132").Value = 35344.656
$ws.Range("I132").Value = 55945.844
$ws.Range("J132").Value = 5235.231
$ws.Range("K132").Value = 167837.532
$ws.Range("L132").Value = 15705.693
$ws.Range("M132").Value = -165307.532
$ws.Range("N132").Value = -20765.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1983.3414
$ws.Range("I132").Value = 1840.174
$ws.Range("J132").Value = 2166.2778
$ws.Range("K132").Value = 5520.522
$ws.Range("L132").Value = 6498.8334
$ws.Range("M132").Value = -2990.522
$ws.Range("N132").Value = -11558.8334
$ws.Range("H136").Value = 608897.9399999999
$ws.Range("I136").Value = 927224.4
$ws.Range("J136").Value = 2561.8572
$ws.Range("K136").Value = 2781673.2
$ws.Range("L136").Value = 7685.571599999999
$ws.Range("M136").Value = -2779123.2
$ws.Range("N136").Value = -12785.5716
